$d = $word.ActiveDocument

# Locate the paragraph that holds the "Version X.Y.Z, YYYY-MM-DD" line.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Version*") {
        $target = $p
    }
}

$pStart = $target.Range.Start

# Within that paragraph only, replace the specific digits that changed:
#   11.02.03, 2015-05-27  ->  11.03.05, 2015-06-14
# "Version 11.02.03, 2015-05-27"
#  0123456789012345678901234567890
#            1111111111222222222

# "2" (index 12, the 2 in ".02") -> "3"
$r = $d.Range($pStart + 12, $pStart + 13)
$r.Text = "3"

# "3" (index 15, the 3 in ".03") -> "5"
$r = $d.Range($pStart + 15, $pStart + 16)
$r.Text = "5"

# "5" (index 24, the 5 in "-05-") -> "6"
$r = $d.Range($pStart + 24, $pStart + 25)
$r.Text = "6"

# "27" (index 26-27, the day) -> "14"
$r = $d.Range($pStart + 26, $pStart + 28)
$r.Text = "14"
